# 1. Hardcoded values moved to Properties:
#    "gyanedel" (Username/Password on TC_010) -> "harikdel"
#    Mobilenumber on TC_010 -> "824860691" (new value, previously shared the
#    generic "8248606917" used by the other rows)
$wb = $excel.ActiveWorkbook

$wsTestSuite = $wb.Worksheets.Item("TestSuite")
$wsLogin     = $wb.Worksheets.Item("TC001_VerifyLogin")
$wsImgReq    = $wb.Worksheets.Item("TC001_VerifyImagerequired")

$wsLogin.Range("B11").Value = "harikdel"
$wsLogin.Range("C11").Value = "harikdel"
$wsLogin.Range("E11").Value = "824860691"

# 2./3. Screenshot/report folder + mail changes don't affect this workbook's
# data, only the view/selection state recorded when the workbook was last
# saved changed (reviewer clicked around the sheets). Reproduce that by
# activating each sheet and selecting the recorded cell, finishing on the
# sheet that ends up as the active tab.

[void]$wsTestSuite.Activate()
[void]$wsTestSuite.Range("I9").Select()

[void]$wsImgReq.Activate()
[void]$wsImgReq.Range("G14").Select()

[void]$wsLogin.Activate()
[void]$wsLogin.Range("H12").Select()
